# Auto-generated edit script: updates crypto price/volume table per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.314.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.92%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.490.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.29%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.53%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.85%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.489.83"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.28%  "

# Row 9
$ws.Range("E9").Value = "  +1.51%  "

# Row 10
$ws.Range("E10").Value = "  +0.75%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.124"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.42%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.436"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.99%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.092.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.38%  "

# Row 14
$ws.Range("E14").Value = "  +0.50%  "

# Row 15
$ws.Range("E15").Value = "  +2.86%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000178"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.10%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.269.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.72%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.496.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.00%  "

# Row 19
$ws.Range("E19").Value = "  +2.44%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.95%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.70%  "

# Row 22
$ws.Range("E22").Value = "  +3.18%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.55%  "

# Row 24
$ws.Range("E24").Value = "  -0.13%  "

# Row 25
$ws.Range("E25").Value = "  +1.77%  "

# Row 26
$ws.Range("E26").Value = "  +5.38%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.29%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.180"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.41%  "

# Row 29
$ws.Range("E29").Value = "  -0.38%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.28%  "

# Row 31
$ws.Range("E31").Value = "  +4.52%  "

# Row 32
$ws.Range("E32").Value = "  +4.52%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.63%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.47%  "

# Row 36
$ws.Range("E36").Value = "  +1.06%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.08%  "

# Row 38
$ws.Range("E38").Value = "  +9.40%  "

# Row 39
$ws.Range("E39").Value = "  +4.72%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0747"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.55%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.65%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.76%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.29%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.808.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.05%  "

# Row 45
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.30%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.63%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0312"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.27%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.93%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "349.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.76%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.36%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "32.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.24%  "
